$wb = $excel.ActiveWorkbook

# Sheet "Hoja1": update the conversion text in A1
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.33 = 8779.46 pesos`n✅ 8779.46 pesos = 2.32 = 958.44 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas": update rate values
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 428.5
$ws2.Range("O10").Value = 3762
$ws2.Range("N12").Value = 3784.99
$ws2.Range("O12").Value = 413.2
